$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "price" column (C) data over to column E,
# and restore "unit" (currently in D) back into D (unchanged).
$ws.Range("E1").Value2 = $ws.Range("C1").Value2
$ws.Range("E2").Value2 = $ws.Range("C2").Value2

# Insert the new "desc" column content into C
$ws.Range("C1").Value2 = "desc"
$ws.Range("C2").Value2 = "about the service"

# Add the new "note" column in F
$ws.Range("F1").Value2 = "note"
$ws.Range("F2").Value2 = "ملاحظات"

# Column widths (closest achievable values; runtime quantizes to 1/6 char units)
$ws.Columns.Item(3).ColumnWidth = 47/3
$ws.Columns.Item(5).ColumnWidth = 41/3
$ws.Columns.Item(6).ColumnWidth = 83/6

# Selection as recorded in the saved file
$ws.Range("E10").Select()
